$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('Bitcoin','https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc','26.480.51','  +0.63%  ')
    ,@('Ethereum','https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth','1.727.95','  +0.66%  ')
    ,@('TetherUSD','https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt','0.9991','  +0.08%  ')
    ,@('BNB','https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb','245.19','  +1.93%  ')
    ,@('USDC','https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc','0.9995','  +0.02%  ')
    ,@('XRP','https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp','0.4792','  +1.41%  ')
    ,@('Cardano','https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada','0.2686','  +1.64%  ')
    ,@('Dogecoin','https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge','0.06224','  -0.18%  ')
    ,@('WrappedEther','https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth','1.728.41','  +0.81%  ')
    ,@('TRON','https://coinranking.com/coin/qUhEFk1I61atv+tron-trx','0.07135','  +0.69%  ')
    ,@('Solana','https://coinranking.com/coin/zNZHO_Sjf+solana-sol','15.72','  +2.83%  ')
    ,@('Polygon','https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic','0.6192','  +4.42%  ')
    ,@('Polkadot','https://coinranking.com/coin/25W7FG7om+polkadot-dot','4.532','  +2.36%  ')
    ,@('Litecoin','https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc','77.25','  +1.39%  ')
    ,@('Dai','https://coinranking.com/coin/MoTuySvg7+dai-dai','0.9996','  -0.03%  ')
    ,@('WrappedBTC','https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc','26.490.36','  +0.75%  ')
    ,@('BinanceUSD','https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd','0.9998','  +0.09%  ')
    ,@('ShibaInu','https://coinranking.com/coin/xz24e0BjL+shibainu-shib','0.000006954','  +2.20%  ')
    ,@('Avalanche','https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax','11.71','  +0.92%  ')
    ,@('WrappedliquidstakedEther2.0','https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth','1.951.51','  +1.57%  ')
    ,@('Uniswap','https://coinranking.com/coin/_H5FVG9iW+uniswap-uni','4.539','  -0.64%  ')
    ,@('Cosmos','https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom','8.943','  +1.04%  ')
    ,@('Chainlink','https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link','5.300','  -0.91%  ')
    ,@('Monero','https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr','136.36','  +0.62%  ')
    ,@('EthereumClassic','https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc','15.37','  +1.09%  ')
    ,@('LidoDAOToken','https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo','1.803','  +2.26%  ')
    ,@('Toncoin','https://coinranking.com/coin/67YlI0K1b+toncoin-ton','1.406','  +0.23%  ')
    ,@('BitcoinCash','https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch','106.75','  +0.31%  ')
    ,@('InternetComputer(DFINITY)','https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp','3.982','  -1.43%  ')
    ,@('Stellar','https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm','0.08031','  +3.42%  ')
    ,@('Filecoin','https://coinranking.com/coin/ymQub4fuB+filecoin-fil','3.734','  +0.93%  ')
    ,@('Hedera','https://coinranking.com/coin/jad286TjB+hedera-hbar','0.04561','  +3.32%  ')
    ,@('Frax','https://coinranking.com/coin/KfWtaeV1W+frax-frax','0.9993','  +0.02%  ')
    ,@('HuobiToken','https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht','2.617','  +0.23%  ')
    ,@('ImmutableX','https://coinranking.com/coin/Z96jIvLU7+immutablex-imx','0.6377','  +2.55%  ')
    ,@('ARBITRUM','https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb','0.9906','  +2.03%  ')
    ,@('TrustWalletToken','https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt','0.9342','  +1.81%  ')
    ,@('RenderToken','https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr','2.088','  +9.03%  ')
    ,@('MXToken','https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx','2.412','  +0.26%  ')
    ,@('Quant','https://coinranking.com/coin/bauj_21eYVwso+quant-qnt','106.02','  -5.36%  ')
    ,@('PaxDollar','https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp','1.006','  +0.55%  ')
    ,@('FraxShare','https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs','5.725','  +11.05%  ')
    ,@('VeChain','https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet','0.01499','  +2.02%  ')
    ,@('TheSandbox','https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand','0.3915','  +2.45%  ')
    ,@('Aptos','https://coinranking.com/coin/HGYj5JCv5+aptos-apt','6.986','  +11.80%  ')
    ,@('Algorand','https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo','0.1192','  +4.40%  ')
    ,@('Cronos','https://coinranking.com/coin/65PHZTpmE55b+cronos-cro','0.05320','  +0.46%  ')
    ,@('Elrond','https://coinranking.com/coin/omwkOTglq+elrond-egld','31.06','  +1.04%  ')
    ,@('EnergySwap','https://coinranking.com/coin/SbWqqTui-+energyswap-ens','7.876','  +3.00%  ')
    ,@('NEARProtocol','https://coinranking.com/coin/DCrsaMv68+nearprotocol-near','1.269','  +3.62%  ')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
}
